# Insert a new price-record row at row 97 of the "Choclo" price sheet,
# shifting the existing rows 97-152 down by one (to 98-153), and
# populate the new row 97 with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 152 -> 153, 151 -> 152, ..., 97 -> 98 (work bottom-up so
# we never overwrite data before it has been copied down).
for ($r = 152; $r -ge 97; $r--) {
    $srcRow = $r
    $dstRow = $r + 1
    $src = $ws.Range("A" + $srcRow + ":R" + $srcRow)
    $dst = $ws.Range("A" + $dstRow + ":R" + $dstRow)
    $dst.Value = $src.Value()
}

# Populate the new row 97 with the new record.
$ws.Range("A97").Value = 4
$ws.Range("B97").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C97").Value = "Los Lagos"
$ws.Range("D97").Value = 44523
$ws.Range("E97").Value = 10
$ws.Range("F97").Value = 100112024
$ws.Range("G97").Value = "Choclo"
$ws.Range("H97").Value = "Dulce o Americano"
$ws.Range("I97").Value = "Primera"
$ws.Range("J97").Value = 250
$ws.Range("K97").Value = 25000
$ws.Range("L97").Value = 25000
$ws.Range("M97").Value = 25000
$ws.Range("N97").Value = "$/malla 70 unidades"
$ws.Range("O97").Value = "Región de Arica y Parinacota"
$ws.Range("P97").Value = 357
$ws.Range("Q97").Value = 70
$ws.Range("R97").Value = "Hortaliza"

# Make sure every date cell in the shifted block keeps the same date
# number format as the rest of the column (the newly created row at
# the bottom can otherwise pick up a different auto-detected format).
$fmt = $ws.Range("D2").NumberFormat
$ws.Range("D97:D153").NumberFormat = $fmt
